$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Efnb1"
$ws.Range("C2").Value2 = "Epha4"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 10.31211433333333
$ws.Range("H2").Value2 = 30.936343
$ws.Range("I2").Value2 = 0.633340936097251
$ws.Range("J2").Value2 = 0.633340936097251
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 6.708176333333333
$ws.Range("N2").Value2 = 20.124529
$ws.Range("O2").Value2 = 0.4356329228871633
$ws.Range("P2").Value2 = 0.4356329228871633
$ws.Range("Q2").Value2 = 69.17548131749412
$ws.Range("R2").Value2 = 622.579331857447
$ws.Range("S2").Value2 = 0.2759041631761375
$ws.Range("T2").Value2 = 0.2759041631761375

$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Efnb1"
$ws.Range("C3").Value2 = "Epha4"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 10.31211433333333
$ws.Range("H3").Value2 = 30.936343
$ws.Range("I3").Value2 = 0.633340936097251
$ws.Range("J3").Value2 = 0.633340936097251
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 6.789877333333333
$ws.Range("N3").Value2 = 20.369632
$ws.Range("O3").Value2 = 0.4409386339573907
$ws.Range("P3").Value2 = 0.4409386339573907
$ws.Range("Q3").Value2 = 70.01799137064178
$ws.Range("R3").Value2 = 630.161922335776
$ws.Range("S3").Value2 = 0.279264487192017
$ws.Range("T3").Value2 = 0.279264487192017

$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Efnb1"
$ws.Range("C4").Value2 = "Epha4"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 10.31211433333333
$ws.Range("H4").Value2 = 30.936343
$ws.Range("I4").Value2 = 0.633340936097251
$ws.Range("J4").Value2 = 0.633340936097251
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 1.900636333333334
$ws.Range("N4").Value2 = 5.701909000000001
$ws.Range("O4").Value2 = 0.1234284431554459
$ws.Range("P4").Value2 = 0.1234284431554459
$ws.Range("Q4").Value2 = 19.59957917542078
$ws.Range("R4").Value2 = 176.396212578787
$ws.Range("S4").Value2 = 0.07817228572909644
$ws.Range("T4").Value2 = 0.07817228572909644

$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Efnb1"
$ws.Range("C5").Value2 = "Epha4"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 4.103438
$ws.Range("H5").Value2 = 12.310314
$ws.Range("I5").Value2 = 0.2520215719230645
$ws.Range("J5").Value2 = 0.2520215719230645
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 6.708176333333333
$ws.Range("N5").Value2 = 20.124529
$ws.Range("O5").Value2 = 0.4356329228871633
$ws.Range("P5").Value2 = 0.4356329228871633
$ws.Range("Q5").Value2 = 27.52658567690067
$ws.Range("R5").Value2 = 247.739271092106
$ws.Range("S5").Value2 = 0.109788894007462
$ws.Range("T5").Value2 = 0.109788894007462

$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Efnb1"
$ws.Range("C6").Value2 = "Epha4"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 4.103438
$ws.Range("H6").Value2 = 12.310314
$ws.Range("I6").Value2 = 0.2520215719230645
$ws.Range("J6").Value2 = 0.2520215719230645
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 6.789877333333333
$ws.Range("N6").Value2 = 20.369632
$ws.Range("O6").Value2 = 0.4409386339573907
$ws.Range("P6").Value2 = 0.4409386339573907
$ws.Range("Q6").Value2 = 27.86184066493866
$ws.Range("R6").Value2 = 250.7565659844479
$ws.Range("S6").Value2 = 0.1111260476515504
$ws.Range("T6").Value2 = 0.1111260476515503

$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Efnb1"
$ws.Range("C7").Value2 = "Epha4"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 4.103438
$ws.Range("H7").Value2 = 12.310314
$ws.Range("I7").Value2 = 0.2520215719230645
$ws.Range("J7").Value2 = 0.2520215719230645
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 1.900636333333334
$ws.Range("N7").Value2 = 5.701909000000001
$ws.Range("O7").Value2 = 0.1234284431554459
$ws.Range("P7").Value2 = 0.1234284431554459
$ws.Range("Q7").Value2 = 7.799143354380667
$ws.Range("R7").Value2 = 70.192290189426
$ws.Range("S7").Value2 = 0.03110663026405209
$ws.Range("T7").Value2 = 0.03110663026405209

$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Efnb1"
$ws.Range("C8").Value2 = "Epha4"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 1.866538
$ws.Range("H8").Value2 = 5.599614
$ws.Range("I8").Value2 = 0.1146374919796846
$ws.Range("J8").Value2 = 0.1146374919796846
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 0.6666666666666666
$ws.Range("M8").Value2 = 6.708176333333333
$ws.Range("N8").Value2 = 20.124529
$ws.Range("O8").Value2 = 0.4356329228871633
$ws.Range("P8").Value2 = 0.4356329228871633
$ws.Range("Q8").Value2 = 12.52106603686733
$ws.Range("R8").Value2 = 112.689594331806
$ws.Range("S8").Value2 = 0.04993986570356374
$ws.Range("T8").Value2 = 0.04993986570356373

$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Efnb1"
$ws.Range("C9").Value2 = "Epha4"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 1.866538
$ws.Range("H9").Value2 = 5.599614
$ws.Range("I9").Value2 = 0.1146374919796846
$ws.Range("J9").Value2 = 0.1146374919796846
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 6.789877333333333
$ws.Range("N9").Value2 = 20.369632
$ws.Range("O9").Value2 = 0.4409386339573907
$ws.Range("P9").Value2 = 0.4409386339573907
$ws.Range("Q9").Value2 = 12.67356405800533
$ws.Range("R9").Value2 = 114.062076522048
$ws.Range("S9").Value2 = 0.05054809911382346
$ws.Range("T9").Value2 = 0.05054809911382346

$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Efnb1"
$ws.Range("C10").Value2 = "Epha4"
$ws.Range("D10").Value2 = "sCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 1.866538
$ws.Range("H10").Value2 = 5.599614
$ws.Range("I10").Value2 = 0.1146374919796846
$ws.Range("J10").Value2 = 0.1146374919796846
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 1.900636333333334
$ws.Range("N10").Value2 = 5.701909000000001
$ws.Range("O10").Value2 = 0.1234284431554459
$ws.Range("P10").Value2 = 0.1234284431554459
$ws.Range("Q10").Value2 = 3.547609940347334
$ws.Range("R10").Value2 = 31.928489463126
$ws.Range("S10").Value2 = 0.01414952716229739
$ws.Range("T10").Value2 = 0.01414952716229739
